# Update the "Förändrad" (Changed) date column (C) for rows 2-18
# from 2023-09-15 (45184) to 2023-09-16 (45185).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45185
    }
}
